$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9, 10, 11 need to be cyclically shifted:
#   old row 9  (byte_pool_create)          -> new row 10
#   old row 10 (thread_preemption_change)  -> new row 11
#   old row 11 (timer_create)              -> new row 9
#
# Capture the original values first (as text), then write them to their
# new positions.

$row9  = @($ws.Cells.Item(9,1).Text,  $ws.Cells.Item(9,2).Text,  $ws.Cells.Item(9,3).Text,  $ws.Cells.Item(9,4).Text,  $ws.Cells.Item(9,5).Text)
$row10 = @($ws.Cells.Item(10,1).Text, $ws.Cells.Item(10,2).Text, $ws.Cells.Item(10,3).Text, $ws.Cells.Item(10,4).Text, $ws.Cells.Item(10,5).Text)
$row11 = @($ws.Cells.Item(11,1).Text, $ws.Cells.Item(11,2).Text, $ws.Cells.Item(11,3).Text, $ws.Cells.Item(11,4).Text, $ws.Cells.Item(11,5).Text)

# New row 9 = old row 11 (timer_create)
$ws.Cells.Item(9,1).Value = $row11[0]
$ws.Cells.Item(9,2).Value = $row11[1]
$ws.Cells.Item(9,3).Value = $row11[2]
$ws.Cells.Item(9,4).Value = [double]$row11[3]
$ws.Cells.Item(9,5).Value = $row11[4]

# New row 10 = old row 9 (byte_pool_create)
$ws.Cells.Item(10,1).Value = $row9[0]
$ws.Cells.Item(10,2).Value = $row9[1]
$ws.Cells.Item(10,3).Value = $row9[2]
$ws.Cells.Item(10,4).Value = [double]$row9[3]
$ws.Cells.Item(10,5).Value = $row9[4]

# New row 11 = old row 10 (thread_preemption_change).
# "60,208" in column C looks like a thousands-grouped number to Excel's
# parser, so format the cell as text first, then restore its original
# (unstyled) appearance afterwards.
$ws.Cells.Item(11,1).Value = $row10[0]
$ws.Cells.Item(11,2).Value = $row10[1]
$ws.Cells.Item(11,3).NumberFormat = "@"
$ws.Cells.Item(11,3).Value = $row10[2]
$ws.Cells.Item(11,3).Style = "Normal"
$ws.Cells.Item(11,4).Value = [double]$row10[3]
$ws.Cells.Item(11,5).Value = $row10[4]
